# Generate Report for handoff
# - Rename source file d7905a2c-...md -> d9d23f4b-...md (and the derived .xlf names / hash)
# - Update handoff datetimes
# - Remove the "Handoff transform failed" row (a96a201b-....md), shifting the
#   ".localization-config" row up to take its place.

$wb = $excel.ActiveWorkbook

$newGuid = "d9d23f4b-30d0-4a4d-9f48-2dc227f8c989"
$newHash = "4d3c858dcf09e105075ad339a34a37c2101b4a87"

# --- Overview sheet ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Rows.Item(3).Delete()

# --- zh-cn sheet ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("C2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-01-14 03:33:34"
$ws2.Rows.Item(3).Delete()

# --- de-de sheet ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("C2").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("D2").Value = "2016-01-14 03:33:45"
$ws3.Rows.Item(3).Delete()
